$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Carel_Partitions")
Write-Host $ws.Name
